$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph 1 : the "{m:comment Style Titre 1}" field is rewritten from a
# real Word field (fldChar/instrText) into plain literal text runs that
# spell out the M2Doc token syntax directly ("{m:comment " / "Style Titre 1"
# / "}"), keeping the same run formatting (lang, and the Titre1Car character
# style on the middle run).
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(1)
$p1Start = $p1.Range.Start
$p1End = $p1.Range.End
$r1 = $d.Range($p1Start, $p1End - 1)

$xml1 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="007655F8" w:rsidRDefault="007655F8" w:rsidP="00F5495F">
<w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">{m:comment </w:t></w:r>
<w:r w:rsidRPr="007655F8"><w:t>Style Titre 1</w:t></w:r>
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r1.InsertXML($xml1)

# Re-apply the "Titre 1 Car" character style on the "Style Titre 1" run
# (InsertXML drops unresolved w:rStyle references, so set it back through
# the object model, exactly like the original instrText run had it).
$p1b = $d.Paragraphs.Item(1)
$prefixLen = "{m:comment ".Length
$midLen = "Style Titre 1".Length
$styleStart = $p1b.Range.Start + $prefixLen
$styleEnd = $styleStart + $midLen
$midRange = $d.Range($styleStart, $styleEnd)
$midRange.Style = "Titre 1 Car"

# ---------------------------------------------------------------------------
# Paragraph 3 : the "{m:'Section 1'.asParagraph().setStyle('Titre1')}"
# field is rewritten the same way, one literal-text run per former
# instrText run (the leading/trailing space runs merge into the "{" / "}"
# delimiter runs, same as the begin/end fldChar runs they used to carry).
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(3)
$p3Start = $p3.Range.Start
$p3End = $p3.Range.End
$r3 = $d.Range($p3Start, $p3End - 1)

$xml3 = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">
<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>{</w:t></w:r>
<w:r w:rsidR="00DE6D5A" w:rsidRPr="009E41B4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>m</w:t></w:r>
<w:r w:rsidR="00681A3B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">:'Section </w:t></w:r>
<w:r w:rsidR="002F3821"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r>
<w:r w:rsidR="00681A3B"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>'</w:t></w:r>
<w:r w:rsidR="00174936"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.asParagraph()</w:t></w:r>
<w:r w:rsidR="009E41B4" w:rsidRPr="009E41B4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>.</w:t></w:r>
<w:r w:rsidR="00174936"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>set</w:t></w:r>
<w:r w:rsidR="007655F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Style</w:t></w:r>
<w:r w:rsidR="009E41B4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r>
<w:r w:rsidR="007655F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>'</w:t></w:r>
<w:r w:rsidR="00B22C70"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Titre</w:t></w:r>
<w:r w:rsidR="002775EA"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>1</w:t></w:r>
<w:r w:rsidR="007655F8"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>'</w:t></w:r>
<w:r w:rsidR="009E41B4"><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>)</w:t></w:r>
<w:r><w:t xml:space="preserve">}</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$r3.InsertXML($xml3)

Write-Output "done"
